$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.472.23'
$ws.Range("E2").Value = '  +0.94%  '

$ws.Range("D3").Value = '1.918.54'
$ws.Range("E3").Value = '  +1.53%  '

$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  +0.78%  '

$ws.Range("D5").Value = '325.45'
$ws.Range("E5").Value = '  +0.99%  '

$ws.Range("E6").Value = '  +0.67%  '

$ws.Range("D7").Value = '0.4827'
$ws.Range("E7").Value = '  +2.67%  '

$ws.Range("D8").Value = '0.4078'
$ws.Range("E8").Value = '  +1.34%  '

$ws.Range("E9").Value = '  +2.11%  '

$ws.Range("D10").Value = '1.021'
$ws.Range("E10").Value = '  +2.66%  '

$ws.Range("D11").Value = '23.49'
$ws.Range("E11").Value = '  +3.51%  '

$ws.Range("D12").Value = '1.906.18'
$ws.Range("E12").Value = '  +1.39%  '

$ws.Range("D13").Value = '6.036'
$ws.Range("E13").Value = '  +2.17%  '

$ws.Range("D14").Value = '7.218'
$ws.Range("E14").Value = '  +2.76%  '

$ws.Range("D15").Value = '91.25'
$ws.Range("E15").Value = '  +1.61%  '

$ws.Range("D16").Value = '0.06789'
$ws.Range("E16").Value = '  +2.36%  '

$ws.Range("D17").Value = '1.008'
$ws.Range("E17").Value = '  +0.68%  '

$ws.Range("D18").Value = '0.00001038'
$ws.Range("E18").Value = '  +1.62%  '

$ws.Range("E19").Value = '  +1.79%  '

$ws.Range("D21").Value = '29.503.76'
$ws.Range("E21").Value = '  +1.04%  '

$ws.Range("D22").Value = '5.643'
$ws.Range("E22").Value = '  +2.80%  '

$ws.Range("D23").Value = '11.76'
$ws.Range("E23").Value = '  +1.01%  '

$ws.Range("D24").Value = '2.193'
$ws.Range("E24").Value = '  +0.80%  '

$ws.Range("D25").Value = '2.142.49'
$ws.Range("E25").Value = '  +1.29%  '

$ws.Range("D26").Value = '6.751'
$ws.Range("E26").Value = '  +11.99%  '

$ws.Range("D27").Value = '156.82'
$ws.Range("E27").Value = '  +1.02%  '

$ws.Range("D28").Value = '20.04'
$ws.Range("E28").Value = '  +1.97%  '

$ws.Range("D29").Value = '2.113'
$ws.Range("E29").Value = '  +1.53%  '

$ws.Range("D30").Value = '120.23'
$ws.Range("E30").Value = '  +0.99%  '

$ws.Range("D31").Value = '1.024'
$ws.Range("E31").Value = '  -0.37%  '

$ws.Range("D32").Value = '0.09565'
$ws.Range("E32").Value = '  +1.60%  '

$ws.Range("D33").Value = '5.522'
$ws.Range("E33").Value = '  +3.06%  '

$ws.Range("D34").Value = '3.564'
$ws.Range("E34").Value = '  +0.68%  '

$ws.Range("D35").Value = '1.382'
$ws.Range("E35").Value = '  +0.08%  '

$ws.Range("D36").Value = "'0.02280"
$ws.Range("E36").Value = '  +2.13%  '

$ws.Range("D37").Value = '0.06134'

$ws.Range("D38").Value = '1.183'
$ws.Range("E38").Value = '  +1.09%  '

$ws.Range("D39").Value = '0.5983'
$ws.Range("E39").Value = '  +2.90%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '8.021'
$ws.Range("E40").Value = '  +0.27%  '

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '10.79'
$ws.Range("E41").Value = '  +7.52%  '

$ws.Range("D42").Value = '0.1854'
$ws.Range("E42").Value = '  +1.36%  '

$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = '1.282'
$ws.Range("E43").Value = '  +0.67%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '2.392'
$ws.Range("E44").Value = '  -2.46%  '

$ws.Range("D45").Value = "'0.07610"
$ws.Range("E45").Value = '  -1.29%  '

$ws.Range("D46").Value = "'12.40"
$ws.Range("E46").Value = '  +1.87%  '

$ws.Range("D47").Value = '0.5573'
$ws.Range("E47").Value = '  +1.91%  '

$ws.Range("D48").Value = '1.956'
$ws.Range("E48").Value = '  +2.86%  '

$ws.Range("D49").Value = '117.51'
$ws.Range("E49").Value = '  +3.82%  '

$ws.Range("E50").Value = '  +4.19%  '

$ws.Range("D51").Value = '72.64'
$ws.Range("E51").Value = '  +2.20%  '
